# Fruta / hortaliza, semanal
#
# Insert two new weekly observations at the top of the "Palta" (Hass) data
# block for Terminal Hortofrutícola Agro Chillán (rows 683:684), pushing the
# existing historical rows down by two positions. Then populate the two new
# rows with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data (previously rows 683:738) down by two rows,
# carrying formatting (e.g. the date style on column D) along with it.
$ws.Range("A683:T684").Insert()

# New row 683: Primera, $/bandeja 10 kilos, origin Perú
$ws.Cells.Item(683, 1).Value = 7
$ws.Cells.Item(683, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(683, 3).Value = "Ñuble"
$ws.Cells.Item(683, 4).Value = 45013
$ws.Cells.Item(683, 5).Value = 16
$ws.Cells.Item(683, 6).Value = "Fruta"
$ws.Cells.Item(683, 7).Value = 100106
$ws.Cells.Item(683, 8).Value = "Oleaginosos"
$ws.Cells.Item(683, 9).Value = 100106002
$ws.Cells.Item(683, 10).Value = "Palta"
$ws.Cells.Item(683, 11).Value = "Hass"
$ws.Cells.Item(683, 12).Value = "Primera"
$ws.Cells.Item(683, 13).Value = 160
$ws.Cells.Item(683, 14).Value = 36000
$ws.Cells.Item(683, 15).Value = 36000
$ws.Cells.Item(683, 16).Value = 36000
$ws.Cells.Item(683, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(683, 18).Value = "Perú"
$ws.Cells.Item(683, 19).Value = 3600
$ws.Cells.Item(683, 20).Value = 10

# New row 684: Primera, $/kilo (en caja de 17 kilos), origin Provincia de Quillota
$ws.Cells.Item(684, 1).Value = 7
$ws.Cells.Item(684, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(684, 3).Value = "Ñuble"
$ws.Cells.Item(684, 4).Value = 45013
$ws.Cells.Item(684, 5).Value = 16
$ws.Cells.Item(684, 6).Value = "Fruta"
$ws.Cells.Item(684, 7).Value = 100106
$ws.Cells.Item(684, 8).Value = "Oleaginosos"
$ws.Cells.Item(684, 9).Value = 100106002
$ws.Cells.Item(684, 10).Value = "Palta"
$ws.Cells.Item(684, 11).Value = "Hass"
$ws.Cells.Item(684, 12).Value = "Primera"
$ws.Cells.Item(684, 13).Value = 130
$ws.Cells.Item(684, 14).Value = 4700
$ws.Cells.Item(684, 15).Value = 5000
$ws.Cells.Item(684, 16).Value = 4862
$ws.Cells.Item(684, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(684, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(684, 19).Value = 4862
$ws.Cells.Item(684, 20).Value = 1
